$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G (K column), rows 2-31, replacing old Strike# derived values
$gValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 2
    8  = 1
    9  = 2
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 2
    15 = 2
    16 = 2
    17 = 1
    18 = 3
    19 = 2
    20 = 0
    21 = 1
    22 = 4
    23 = 0
    24 = 2
    25 = 1
    26 = 2
    27 = 5
    28 = 1
    29 = 3
    30 = 0
    31 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
